$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 16.71895933333333
$ws.Range("H2").Value = 50.156878
$ws.Range("I2").Value = 0.02912144738161902
$ws.Range("J2").Value = 0.03059269312988411
$ws.Range("M2").Value = 46.17354133333333
$ws.Range("N2").Value = 138.520624
$ws.Range("O2").Value = 0.1154336358852217
$ws.Range("P2").Value = 0.1189208138601986
$ws.Range("Q2").Value = 771.9735598279858
$ws.Range("R2").Value = 6947.762038451871
$ws.Range("S2").Value = 0.003361594553500453
$ws.Range("T2").Value = 0.003638107965181125

# Row 3
$ws.Range("G3").Value = 16.71895933333333
$ws.Range("H3").Value = 50.156878
$ws.Range("I3").Value = 0.02912144738161902
$ws.Range("J3").Value = 0.03059269312988411
$ws.Range("O3").Value = 0.2100327918507284
$ws.Range("P3").Value = 0.2163777511873036
$ws.Range("Q3").Value = 1404.61452818515
$ws.Range("R3").Value = 12641.53075366635
$ws.Range("S3").Value = 0.006116458896295526
$ws.Range("T3").Value = 0.006619578142207596

# Row 4
$ws.Range("G4").Value = 16.71895933333333
$ws.Range("H4").Value = 50.156878
$ws.Range("I4").Value = 0.02912144738161902
$ws.Range("J4").Value = 0.03059269312988411
$ws.Range("M4").Value = 128.0910926666667
$ws.Range("N4").Value = 384.273278
$ws.Range("O4").Value = 0.3202271284388135
$ws.Range("P4").Value = 0.3299009897940278
$ws.Range("Q4").Value = 2141.549769256232
$ws.Range("R4").Value = 19273.94792330608
$ws.Range("S4").Value = 0.009325477470997863
$ws.Range("T4").Value = 0.01009255974401372

# Row 5
$ws.Range("G5").Value = 16.71895933333333
$ws.Range("H5").Value = 50.156878
$ws.Range("I5").Value = 0.02912144738161902
$ws.Range("J5").Value = 0.03059269312988411
$ws.Range("M5").Value = 35.18830149999999
$ws.Range("N5").Value = 70.37660299999999
$ws.Range("O5").Value = 0.08797058803540478
$ws.Range("P5").Value = 0.06041874966919073
$ws.Range("Q5").Value = 588.3117817875723
$ws.Range("R5").Value = 3529.870690725433
$ws.Range("S5").Value = 0.002561830850603124
$ws.Range("T5").Value = 0.001848372267920839

# Row 6
$ws.Range("G6").Value = 16.71895933333333
$ws.Range("H6").Value = 50.156878
$ws.Range("I6").Value = 0.02912144738161902
$ws.Range("J6").Value = 0.03059269312988411
$ws.Range("M6").Value = 106.534543
$ws.Range("N6").Value = 319.603629
$ws.Range("O6").Value = 0.2663358557898317
$ws.Range("P6").Value = 0.2743816954892795
$ws.Range("Q6").Value = 1781.146692012251
$ws.Range("R6").Value = 16030.32022811026
$ws.Range("S6").Value = 0.007756085610222053
$ws.Range("T6").Value = 0.008394075010560834

# Row 7
$ws.Range("I7").Value = 0.2708539632042961
$ws.Range("J7").Value = 0.2845377865576845
$ws.Range("M7").Value = 46.17354133333333
$ws.Range("N7").Value = 138.520624
$ws.Range("O7").Value = 0.1154336358852217
$ws.Range("P7").Value = 0.1189208138601986
$ws.Range("Q7").Value = 7180.003638840913
$ws.Range("R7").Value = 64620.03274956821
$ws.Range("S7").Value = 0.03126565776659396
$ws.Range("T7").Value = 0.03383746515141933

# Row 8
$ws.Range("I8").Value = 0.2708539632042961
$ws.Range("J8").Value = 0.2845377865576845
$ws.Range("O8").Value = 0.2100327918507284
$ws.Range("P8").Value = 0.2163777511873036
$ws.Range("S8").Value = 0.05688821407563276
$ws.Range("T8").Value = 0.06156764638316476

# Row 9
$ws.Range("I9").Value = 0.2708539632042961
$ws.Range("J9").Value = 0.2845377865576845
$ws.Range("M9").Value = 128.0910926666667
$ws.Range("N9").Value = 384.273278
$ws.Range("O9").Value = 0.3202271284388135
$ws.Range("P9").Value = 0.3299009897940278
$ws.Range("Q9").Value = 19918.21473710172
$ws.Range("R9").Value = 179263.9326339154
$ws.Range("S9").Value = 0.08673478686318378
$ws.Range("T9").Value = 0.09386929741918193

# Row 10
$ws.Range("I10").Value = 0.2708539632042961
$ws.Range("J10").Value = 0.2845377865576845
$ws.Range("M10").Value = 35.18830149999999
$ws.Range("N10").Value = 70.37660299999999
$ws.Range("O10").Value = 0.08797058803540478
$ws.Range("P10").Value = 0.06041874966919073
$ws.Range("Q10").Value = 5471.794571499283
$ws.Range("R10").Value = 32830.7674289957
$ws.Range("S10").Value = 0.02382718241480181
$ws.Range("T10").Value = 0.01719141729745436

# Row 11
$ws.Range("I11").Value = 0.2708539632042961
$ws.Range("J11").Value = 0.2845377865576845
$ws.Range("M11").Value = 106.534543
$ws.Range("N11").Value = 319.603629
$ws.Range("O11").Value = 0.2663358557898317
$ws.Range("P11").Value = 0.2743816954892795
$ws.Range("Q11").Value = 16566.16287843723
$ws.Range("R11").Value = 149095.465905935
$ws.Range("S11").Value = 0.07213812208408377
$ws.Range("T11").Value = 0.0780719603064642

# Row 12
$ws.Range("G12").Value = 194.8548433333333
$ws.Range("H12").Value = 584.56453
$ws.Range("I12").Value = 0.3394024086099587
$ws.Range("J12").Value = 0.3565493705749576
$ws.Range("M12").Value = 46.17354133333333
$ws.Range("N12").Value = 138.520624
$ws.Range("O12").Value = 0.1154336358852217
$ws.Range("P12").Value = 0.1189208138601986
$ws.Range("Q12").Value = 8997.138162651858
$ws.Range("R12").Value = 80974.24346386672
$ws.Range("S12").Value = 0.03917845405404922
$ws.Range("T12").Value = 0.0424011413301155

# Row 13
$ws.Range("G13").Value = 194.8548433333333
$ws.Range("H13").Value = 584.56453
$ws.Range("I13").Value = 0.3394024086099587
$ws.Range("J13").Value = 0.3565493705749576
$ws.Range("O13").Value = 0.2100327918507284
$ws.Range("P13").Value = 0.2163777511873036
$ws.Range("Q13").Value = 16370.39353804525
$ws.Range("R13").Value = 147333.5418424073
$ws.Range("S13").Value = 0.07128563544121132
$ws.Range("T13").Value = 0.07714935099225788

# Row 14
$ws.Range("G14").Value = 194.8548433333333
$ws.Range("H14").Value = 584.56453
$ws.Range("I14").Value = 0.3394024086099587
$ws.Range("J14").Value = 0.3565493705749576
$ws.Range("M14").Value = 128.0910926666667
$ws.Range("N14").Value = 384.273278
$ws.Range("O14").Value = 0.3202271284388135
$ws.Range("P14").Value = 0.3299009897940278
$ws.Range("Q14").Value = 24959.16979395882
$ws.Range("R14").Value = 224632.5281456293
$ws.Range("S14").Value = 0.1086858586943839
$ws.Range("T14").Value = 0.1176259902631161

# Row 15
$ws.Range("G15").Value = 194.8548433333333
$ws.Range("H15").Value = 584.56453
$ws.Range("I15").Value = 0.3394024086099587
$ws.Range("J15").Value = 0.3565493705749576
$ws.Range("M15").Value = 35.18830149999999
$ws.Range("N15").Value = 70.37660299999999
$ws.Range("O15").Value = 0.08797058803540478
$ws.Range("P15").Value = 0.06041874966919073
$ws.Range("Q15").Value = 6856.610975948597
$ws.Range("R15").Value = 41139.66585569159
$ws.Range("S15").Value = 0.0298574294660508
$ws.Range("T15").Value = 0.02154226716547588

# Row 16
$ws.Range("G16").Value = 194.8548433333333
$ws.Range("H16").Value = 584.56453
$ws.Range("I16").Value = 0.3394024086099587
$ws.Range("J16").Value = 0.3565493705749576
$ws.Range("M16").Value = 106.534543
$ws.Range("N16").Value = 319.603629
$ws.Range("O16").Value = 0.2663358557898317
$ws.Range("P16").Value = 0.2743816954892795
$ws.Range("Q16").Value = 20758.77168585326
$ws.Range("R16").Value = 186828.9451726794
$ws.Range("S16").Value = 0.09039503095426349
$ws.Range("T16").Value = 0.09783062082399227

# Row 17
$ws.Range("G17").Value = 82.82950199999999
$ws.Range("H17").Value = 165.659004
$ws.Range("I17").Value = 0.1442742299952585
$ws.Range("J17").Value = 0.1010420758958371
$ws.Range("M17").Value = 46.17354133333333
$ws.Range("N17").Value = 138.520624
$ws.Range("O17").Value = 0.1154336358852217
$ws.Range("P17").Value = 0.1189208138601986
$ws.Range("Q17").Value = 3824.531434216416
$ws.Range("R17").Value = 22947.18860529849
$ws.Range("S17").Value = 0.0166540989328934
$ws.Range("T17").Value = 0.01201600589965691

# Row 18
$ws.Range("G18").Value = 82.82950199999999
$ws.Range("H18").Value = 165.659004
$ws.Range("I18").Value = 0.1442742299952585
$ws.Range("J18").Value = 0.1010420758958371
$ws.Range("O18").Value = 0.2100327918507284
$ws.Range("P18").Value = 0.2163777511873036
$ws.Range("Q18").Value = 6958.77772963905
$ws.Range("R18").Value = 41752.6663778343
$ws.Range("S18").Value = 0.03030231931801824
$ws.Range("T18").Value = 0.02186325715763809

# Row 19
$ws.Range("G19").Value = 82.82950199999999
$ws.Range("H19").Value = 165.659004
$ws.Range("I19").Value = 0.1442742299952585
$ws.Range("J19").Value = 0.1010420758958371
$ws.Range("M19").Value = 128.0910926666667
$ws.Range("N19").Value = 384.273278
$ws.Range("O19").Value = 0.3202271284388135
$ws.Range("P19").Value = 0.3299009897940278
$ws.Range("Q19").Value = 10609.72141621585
$ws.Range("R19").Value = 63658.32849729511
$ws.Range("S19").Value = 0.04620052237910256
$ws.Range("T19").Value = 0.03333388084887995

# Row 20
$ws.Range("G20").Value = 82.82950199999999
$ws.Range("H20").Value = 165.659004
$ws.Range("I20").Value = 0.1442742299952585
$ws.Range("J20").Value = 0.1010420758958371
$ws.Range("M20").Value = 35.18830149999999
$ws.Range("N20").Value = 70.37660299999999
$ws.Range("O20").Value = 0.08797058803540478
$ws.Range("P20").Value = 0.06041874966919073
$ws.Range("Q20").Value = 2914.629489470852
$ws.Range("R20").Value = 11658.51795788341
$ws.Range("S20").Value = 0.01269188885103812
$ws.Range("T20").Value = 0.006104835889605954

# Row 21
$ws.Range("G21").Value = 82.82950199999999
$ws.Range("H21").Value = 165.659004
$ws.Range("I21").Value = 0.1442742299952585
$ws.Range("J21").Value = 0.1010420758958371
$ws.Range("M21").Value = 106.534543
$ws.Range("N21").Value = 319.603629
$ws.Range("O21").Value = 0.2663358557898317
$ws.Range("P21").Value = 0.2743816954892795
$ws.Range("Q21").Value = 8824.203142487586
$ws.Range("R21").Value = 52945.21885492552
$ws.Range("S21").Value = 0.03842540051420617
$ws.Range("T21").Value = 0.02772409610005625

# Row 22
$ws.Range("G22").Value = 124.2078576666667
$ws.Range("H22").Value = 372.623573
$ws.Range("I22").Value = 0.2163479508088675
$ws.Range("J22").Value = 0.2272780738416368
$ws.Range("M22").Value = 46.17354133333333
$ws.Range("N22").Value = 138.520624
$ws.Range("O22").Value = 0.1154336358852217
$ws.Range("P22").Value = 0.1189208138601986
$ws.Range("Q22").Value = 5735.116649896616
$ws.Range("R22").Value = 51616.04984906955
$ws.Range("S22").Value = 0.02497383057818468
$ws.Range("T22").Value = 0.02702809351382577

# Row 23
$ws.Range("G23").Value = 124.2078576666667
$ws.Range("H23").Value = 372.623573
$ws.Range("I23").Value = 0.2163479508088675
$ws.Range("J23").Value = 0.2272780738416368
$ws.Range("O23").Value = 0.2100327918507284
$ws.Range("P23").Value = 0.2163777511873036
$ws.Range("Q23").Value = 10435.10890331053
$ws.Range("R23").Value = 93915.98012979471
$ws.Range("S23").Value = 0.0454401641195705
$ws.Range("T23").Value = 0.04917791851203531

# Row 24
$ws.Range("G24").Value = 124.2078576666667
$ws.Range("H24").Value = 372.623573
$ws.Range("I24").Value = 0.2163479508088675
$ws.Range("J24").Value = 0.2272780738416368
$ws.Range("M24").Value = 128.0910926666667
$ws.Range("N24").Value = 384.273278
$ws.Range("O24").Value = 0.3202271284388135
$ws.Range("P24").Value = 0.3299009897940278
$ws.Range("Q24").Value = 15909.92020630914
$ws.Range("R24").Value = 143189.2818567823
$ws.Range("S24").Value = 0.06928048303114534
$ws.Range("T24").Value = 0.07497926151883612

# Row 25
$ws.Range("G25").Value = 124.2078576666667
$ws.Range("H25").Value = 372.623573
$ws.Range("I25").Value = 0.2163479508088675
$ws.Range("J25").Value = 0.2272780738416368
$ws.Range("M25").Value = 35.18830149999999
$ws.Range("N25").Value = 70.37660299999999
$ws.Range("O25").Value = 0.08797058803540478
$ws.Range("P25").Value = 0.06041874966919073
$ws.Range("Q25").Value = 4370.663544243752
$ws.Range("R25").Value = 26223.98126546251
$ws.Range("S25").Value = 0.0190322564529109
$ws.Range("T25").Value = 0.0137318570487337

# Row 26
$ws.Range("G26").Value = 124.2078576666667
$ws.Range("H26").Value = 372.623573
$ws.Range("I26").Value = 0.2163479508088675
$ws.Range("J26").Value = 0.2272780738416368
$ws.Range("M26").Value = 106.534543
$ws.Range("N26").Value = 319.603629
$ws.Range("O26").Value = 0.2663358557898317
$ws.Range("P26").Value = 0.2743816954892795
$ws.Range("Q26").Value = 13232.42735352738
$ws.Range("R26").Value = 16030.32022811026
$ws.Range("S26").Value = 0.05762121662705614
$ws.Range("T26").Value = 0.06236094324820598

